$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.079.73'
$ws.Range("E2").Value = '  +5.76%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.717.60'
$ws.Range("E3").Value = '  +3.61%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '332.74'
$ws.Range("E5").Value = '  +3.78%  '

$ws.Range("E6").Value = '  +0.12%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3685'
$ws.Range("E7").Value = '  +1.26%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '49.29'
$ws.Range("E8").Value = '  +5.19%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3348'
$ws.Range("E9").Value = '  +2.62%  '

$ws.Range("E10").Value = '  +4.54%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07462'
$ws.Range("E11").Value = '  +5.75%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.001'
$ws.Range("E12").Value = '  +0.15%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.291'
$ws.Range("E13").Value = '  +5.08%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.04'
$ws.Range("E14").Value = '  +2.54%  '

$ws.Range("E15").Value = '  +4.52%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.717.26'
$ws.Range("E16").Value = '  +3.40%  '

$ws.Range("E17").Value = '  +2.99%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06632'
$ws.Range("E18").Value = '  +0.04%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '81.84'
$ws.Range("E19").Value = '  +3.73%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.38'
$ws.Range("E21").Value = '  +4.01%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.085'
$ws.Range("E22").Value = '  +2.72%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.01'
$ws.Range("E23").Value = '  +2.93%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '26.019.53'
$ws.Range("E24").Value = '  +5.52%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.474'
$ws.Range("E25").Value = '  +0.45%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.457'
$ws.Range("E26").Value = '  +2.56%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '150.16'
$ws.Range("E27").Value = '  +1.49%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.19'
$ws.Range("E28").Value = '  +3.18%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.319'
$ws.Range("E29").Value = '  +8.88%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.907.41'
$ws.Range("E30").Value = '  +3.47%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '129.27'
$ws.Range("E31").Value = '  +3.01%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.108'
$ws.Range("E32").Value = '  +0.81%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.919'
$ws.Range("E33").Value = '  +1.34%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08531'
$ws.Range("E34").Value = '  +0.69%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.722'
$ws.Range("E35").Value = '  +2.54%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '12.86'
$ws.Range("E36").Value = '  +4.39%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.339'
$ws.Range("E37").Value = '  +2.45%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06215'
$ws.Range("E38").Value = '  +3.08%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02295'
$ws.Range("E39").Value = '  +2.60%  '

$ws.Range("E40").Value = '  +2.70%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.532'
$ws.Range("E41").Value = '  +3.90%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.219'
$ws.Range("E42").Value = '  -4.57%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '14.50'
$ws.Range("E43").Value = '  +13.09%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6163'
$ws.Range("E44").Value = '  +3.93%  '

$ws.Range("E45").Value = '  +0.18%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.833'
$ws.Range("E46").Value = '  -0.66%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5883'
$ws.Range("E47").Value = '  +4.62%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '127.74'
$ws.Range("E48").Value = '  +2.73%  '

$ws.Range("E49").Value = '  +3.26%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07272'
$ws.Range("E50").Value = '  +4.36%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '76.96'
$ws.Range("E51").Value = '  +3.01%  '
